$wb = $excel.ActiveWorkbook

# Rename Sheet3 -> DataSetInteractionPages and make it the active tab
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "DataSetInteractionPages"

# Populate the new data table. Written in this order so the shared
# strings are appended in the same order as the target workbook.
$ws3.Range("A2").Value = "DefaultFunctionality_DragToOppositeCorner_ElementMovedToOppositeCorner"
$ws3.Range("B1").Value = "HorizontalOffset"
$ws3.Range("C1").Value = "VerticalOffset"
$ws3.Range("A1").Value = "Key"
$ws3.Range("B2").Value = 150
$ws3.Range("C2").Value = 140

# Column widths matching the authored layout (closest reachable values).
$ws3.Columns.Item(1).ColumnWidth = 36.5
$ws3.Columns.Item(2).ColumnWidth = 17.8333333333
$ws3.Columns.Item(3).ColumnWidth = 23.5

# Selection on the new sheet, then make it the active/selected tab.
$ws3.Range("C7").Select()
$ws3.Activate()
